$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1): From Currency / To Currency / Exchange Rate / As Of ---
$ws.Range("L1").Value = "From Currency"
$ws.Range("M1").Value = "To Currency"
$ws.Range("N1").Value = "Exchange Rate "
$ws.Range("O1").Value = "As Of"

# --- Row 2 ---
$ws.Range("L2").Value = "USD"
$ws.Range("M2").Value = "INR"
$ws.Range("N2").Value = 80
$ws.Range("G2").Copy()
$ws.Range("O2").PasteSpecial(-4122)
$ws.Range("O2").Value = 44743

# --- Row 3 ---
$ws.Range("L3").Value = "USD"
$ws.Range("M3").Value = "INR"
$ws.Range("N3").Value = 81
$ws.Range("G3").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$ws.Range("O3").Value = 44774

# --- Row 4 ---
$ws.Range("L4").Value = "USD"
$ws.Range("M4").Value = "INR"
$ws.Range("N4").Value = 81
$ws.Range("G4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 44866

$excel.CutCopyMode = $false

# --- Update view: scroll so column C is leftmost, select O4 ---
$ws.Range("O4").Select()
